$d = $word.ActiveDocument

# 1) "mediafile(i.e., playing video)" sentence: merge the spell/grammar-
#    checked runs (proofErr-wrapped "mediafile" / "(" / "i.e., playing
#    video).") back into one continuous run with no proofErr markers.
$d.Content.Find.Execute(
    "The prototype will count the number of fingers displayed in a region of the camera to perform functions for a mediafile(i.e., playing video).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The prototype will count the number of fingers displayed in a region of the camera to perform functions for a mediafile(i.e., playing video).",
    2) | Out-Null

# 2) "Texts, diagrams or pictures would be all fine" - merge the
#    grammar-checked "diagrams" run back into the surrounding text.
$d.Content.Find.Execute(
    "Texts, diagrams or pictures would be all fine",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Texts, diagrams or pictures would be all fine",
    2) | Out-Null

# 3) "Link to a 2-minute Youtube video" - merge the spell-checked
#    "Youtube" run back into the surrounding text.
$d.Content.Find.Execute(
    "Link to a 2-minute Youtube video",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Link to a 2-minute Youtube video",
    2) | Out-Null

# 4) Fix the logical error in the pixel-range equation: "0<=x>=135" should
#    read "0<=x<=135" (the upper bound must use <=, not >=).
$om = $d.OMaths.Item(1)
$mathXml = '<m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Times New Roman" w:hAnsi="Cambria Math" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>0' + [char]0x2264 + 'x</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Times New Roman" w:hAnsi="Cambria Math" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>' + [char]0x2264 + '</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="Times New Roman" w:hAnsi="Cambria Math" w:cs="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t xml:space="preserve">135 , where x is a pixel value. </m:t></m:r></m:oMath>'
$om.Range.InsertXML($mathXml)
